$d = $word.ActiveDocument

$d.Content.Find.Execute("La propbabilidad de pelirrojo es", $true, $false, $false, $false, $false, $true, 1, $false, "La probabilidad de pelirrojo es", 2)

$d.Content.Find.Execute("se puede interperetar como un experimento Bernouilli", $true, $false, $false, $false, $false, $true, 1, $false, "se puede interpretar como un experimento Bernouilli", 2)

$d.Content.Find.Execute("Si consideramos los nacimienro sucesos inpedendientes del mimo experimento", $true, $false, $false, $false, $false, $true, 1, $false, "Si consideramos los nacimientos sucesos independientes del mimo experimento", 2)

$d.Content.Find.Execute("Para el cálculo final hemo utilizado R", $true, $false, $false, $false, $false, $true, 1, $false, "Para el cálculo final hemos utilizado R", 2)

$d.Content.Find.Execute("distribución geimétrica", $true, $false, $false, $false, $false, $true, 1, $false, "distribución geométrica", 2)
